$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; existing rows 7-17 shift down to 8-18
$ws.Rows.Item(7).Insert()

# Copy style (date number format) from the date cell below into the new cell
$ws.Range("D8").Copy()
$ws.Range("D7").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 7 with the new weekly entry (same market/product metadata)
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 44495
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100107
$ws.Cells.Item(7, 8).Value = "Otros"
$ws.Cells.Item(7, 9).Value = 100107002
$ws.Cells.Item(7, 10).Value = "Chirimoya"
$ws.Cells.Item(7, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 100
$ws.Cells.Item(7, 14).Value = 26000
$ws.Cells.Item(7, 15).Value = 27000
$ws.Cells.Item(7, 16).Value = 26500
$ws.Cells.Item(7, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(7, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 19).Value = 2650
$ws.Cells.Item(7, 20).Value = 10
